$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B/C text-only updates (coin name / link swaps) - plain text, safe to assign directly
$ws.Range("B23").Value = 'CoinExToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("B24").Value = 'ZBToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("B42").Value = 'KickToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'

# D/E numeric-looking text values. Force text type via NumberFormat "@" so Excel
# does not auto-convert these into numbers/percentages, then reset the style back
# to Normal afterwards so no stray number-format style is left on the cells.
$numRange = $ws.Range("D2:E47")
$numRange.NumberFormat = "@"

$ws.Range("D2").Value = '279.03'
$ws.Range("E2").Value = '0.42%'
$ws.Range("D3").Value = '27.44'
$ws.Range("E3").Value = '0.32%'
$ws.Range("D4").Value = '4.843'
$ws.Range("E4").Value = '-1.17%'
$ws.Range("D5").Value = '0.06376'
$ws.Range("E5").Value = '0.22%'
$ws.Range("D6").Value = '7.033'
$ws.Range("E6").Value = '0.84%'
$ws.Range("D7").Value = '1.301'
$ws.Range("E7").Value = '2.72%'
$ws.Range("D8").Value = '0.8953'
$ws.Range("E8").Value = '1.51%'
$ws.Range("D9").Value = '0.1523'
$ws.Range("E9").Value = '0.14%'
$ws.Range("D10").Value = '0.05901'
$ws.Range("E10").Value = '15.56%'
$ws.Range("D11").Value = '0.07516'
$ws.Range("E11").Value = '-0.32%'
$ws.Range("D12").Value = '0.02919'
$ws.Range("E12").Value = '-1.61%'
$ws.Range("E13").Value = '-0.13%'
$ws.Range("D14").Value = '0.001567'
$ws.Range("E14").Value = '-0.33%'
$ws.Range("D15").Value = '0.0006388'
$ws.Range("E15").Value = '-0.28%'
$ws.Range("D16").Value = '0.006093'
$ws.Range("E16").Value = '1.24%'
$ws.Range("E17").Value = '0.58%'
$ws.Range("D18").Value = '3.325'
$ws.Range("E18").Value = '0.37%'
$ws.Range("D19").Value = '2.224'
$ws.Range("E19").Value = '-2.13%'
$ws.Range("E21").Value = '1.09%'
$ws.Range("D22").Value = '3.902'
$ws.Range("E22").Value = '-0.28%'
$ws.Range("D23").Value = '0.04417'
$ws.Range("E23").Value = '-0.07%'
$ws.Range("D24").Value = '0.1503'
$ws.Range("E24").Value = '8.90%'
$ws.Range("E25").Value = '0.27%'
$ws.Range("D26").Value = '0.004279'
$ws.Range("E26").Value = '10.52%'
$ws.Range("D28").Value = '0.0001179'
$ws.Range("E28").Value = '-1.69%'
$ws.Range("D29").Value = '0.0001653'
$ws.Range("E29").Value = '-14.62%'
$ws.Range("D40").Value = '0.04073'
$ws.Range("E40").Value = '-1.74%'
$ws.Range("D41").Value = '0.1414'
$ws.Range("E41").Value = '20.01%'
$ws.Range("D42").Value = '0.006627'
$ws.Range("E42").Value = '-3.16%'
$ws.Range("D43").Value = '0.002089'
$ws.Range("E43").Value = '0.95%'
$ws.Range("D44").Value = '0.01100'
$ws.Range("E44").Value = '-2.01%'
$ws.Range("D45").Value = '0.00005529'
$ws.Range("E45").Value = '7.22%'
$ws.Range("D46").Value = '1.561'
$ws.Range("E46").Value = '5.01%'
$ws.Range("D47").Value = '0.01849'
$ws.Range("E47").Value = '-8.67%'

$numRange.Style = "Normal"
